# feat: add 2022-Q1 data
#
# Before: sheets = [ "2021-Q4", "总计" ]
# After:  sheets = [ "2021-Q4", "2022-Q1", "总计" ]
#
# The old "总计" sheet becomes "2022-Q1" (holds the new quarter's fund
# holdings detail) and a brand-new "总计" summary sheet is appended at
# the end, listing both quarters (2022-Q1 first, then 2021-Q4).

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets("总计")

# 1) Duplicate the existing "总计" sheet to the very end of the workbook
#    *before* touching its contents, so the copy still has the old
#    (日期/持有数量(只)/持有市值(亿元)) layout that the new summary sheet needs.
$total.Copy($null, $wb.Worksheets($wb.Worksheets.Count))
$newTotal = $wb.Worksheets($wb.Worksheets.Count)

# 2) Turn the original "总计" sheet into the "2022-Q1" detail sheet.
$total.Name = "2022-Q1"

# 3) Rename the freshly duplicated sheet back to "总计".
$newTotal.Name = "总计"

# ---------------------------------------------------------------------
# Build out "2022-Q1" (fund holdings detail for 601929 as of 2022-Q1)
# ---------------------------------------------------------------------

function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# extend the bold/border header styling from D1 across the new columns
$total.Range("D1").Copy()
$total.Range("E1:H1").PasteSpecial(-4122)

$total.Range("B1").Value = "基金代码"
$total.Range("C1").Value = "基金名称"
$total.Range("D1").Value = "基金规模"
$total.Range("E1").Value = "股票总仓位"
$total.Range("F1").Value = "仓位占比"
$total.Range("G1").Value = "持有市值(亿元)"
$total.Range("H1").Value = "仓位排名"

Set-TextCell $total.Range("B2") "005443"
Set-TextCell $total.Range("C2") "国金量化多策略灵活配置混合"
Set-TextCell $total.Range("D2") "0.51"
Set-TextCell $total.Range("E2") "64.10"
Set-TextCell $total.Range("F2") "0.90"
Set-TextCell $total.Range("G2") "0.0046"
$total.Range("H2").Value = 6

# ---------------------------------------------------------------------
# Build out the new "总计" summary sheet: add a 2022-Q1 row above the
# pre-existing 2021-Q4 row.
# ---------------------------------------------------------------------

$oldB = $newTotal.Range("B2").Value()
$oldC = $newTotal.Range("C2").Value()
$oldD = $newTotal.Range("D2").Value()

# row-index column: the pushed-down 2021-Q4 row advances to index 1
$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = $oldB
$newTotal.Range("C3").Value = $oldC
$newTotal.Range("D3").Value = $oldD

# carry the index-column style (bold/border) from A2 down to A3
$newTotal.Range("A2").Copy()
$newTotal.Range("A3").PasteSpecial(-4122)

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 1
$newTotal.Range("D2").Value = 0
